# Update ePHASORSIM for GMDM profile; fixed GLM tank phasing
# Swap bus name (column A) and angle (column E) between paired rows on the
# "Bus" sheet so that the phase letters line up with the correct angle.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Bus")

# Each entry is a pair of row numbers whose A (name) and E (angle) values
# must be swapped with each other.
$rowPairs = @(
    @(4, 5),
    @(7, 8),
    @(10, 11),
    @(17, 18),
    @(21, 22),
    @(24, 25),
    @(27, 28),
    @(30, 31),
    @(32, 33),
    @(35, 36),
    @(38, 39),
    @(40, 41),
    @(43, 44),
    @(46, 47),
    @(49, 50),
    @(53, 54)
)

foreach ($pair in $rowPairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]

    $nameCell1 = $ws.Cells.Item($r1, 1)
    $nameCell2 = $ws.Cells.Item($r2, 1)
    $angleCell1 = $ws.Cells.Item($r1, 5)
    $angleCell2 = $ws.Cells.Item($r2, 5)

    $name1 = $nameCell1.Value2
    $name2 = $nameCell2.Value2
    $angle1 = $angleCell1.Value2
    $angle2 = $angleCell2.Value2

    $nameCell1.Value2 = $name2
    $nameCell2.Value2 = $name1
    $angleCell1.Value2 = $angle2
    $angleCell2.Value2 = $angle1
}
